$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $val) {
    $r = $ws.Range($rangeAddr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range("D2").Value = '60.857.11'
$ws.Range("E2").Value = '  -1.23%  '
$ws.Range("D3").Value = '3.370.26'
$ws.Range("E3").Value = '  -0.58%  '
$ws.Range("E4").Value = '  +0.05%  '
Set-TextValue "D5" '569.60'
$ws.Range("E5").Value = '  -1.36%  '
Set-TextValue "D6" '136.28'
$ws.Range("E6").Value = '  -0.50%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '3.368.16'
$ws.Range("E8").Value = '  -0.62%  '
$ws.Range("E9").Value = '  -1.49%  '
Set-TextValue "D10" '7.59'
$ws.Range("E10").Value = '  +1.24%  '
$ws.Range("E11").Value = '  -3.56%  '
Set-TextValue "D12" '0.378'
$ws.Range("E12").Value = '  -2.86%  '
$ws.Range("D13").Value = '3.945.30'
$ws.Range("E13").Value = '  -0.54%  '
$ws.Range("E14").Value = '  +0.44%  '
Set-TextValue "D15" '25.57'
$ws.Range("E15").Value = '  -0.17%  '
$ws.Range("D16").Value = '3.373.65'
$ws.Range("E16").Value = '  -0.54%  '
$ws.Range("E17").Value = '  -3.60%  '
$ws.Range("D18").Value = '61.062.43'
$ws.Range("E18").Value = '  -1.09%  '
Set-TextValue "D19" '13.73'
$ws.Range("E19").Value = '  -3.11%  '
Set-TextValue "D20" '5.71'
$ws.Range("E20").Value = '  -1.52%  '
Set-TextValue "D21" '9.30'
$ws.Range("E21").Value = '  -2.02%  '
Set-TextValue "D22" '372.59'
$ws.Range("E22").Value = '  -1.39%  '
$ws.Range("D23").Value = '3.514.55'
$ws.Range("E23").Value = '  -0.24%  '
Set-TextValue "D24" '0.545'
$ws.Range("E24").Value = '  -2.86%  '
Set-TextValue "D25" '0.998'
$ws.Range("E25").Value = '  -0.29%  '
Set-TextValue "D26" '70.87'
$ws.Range("E26").Value = '  -0.40%  '
Set-TextValue "D27" '0.0000122'
$ws.Range("E27").Value = '  -2.13%  '
$ws.Range("E28").Value = '  +11.23%  '
$ws.Range("E29").Value = '  -6.67%  '
Set-TextValue "D30" '1.00'
$ws.Range("E30").Value = '  +0.00%  '
Set-TextValue "D31" '7.31'
$ws.Range("E31").Value = '  -3.94%  '
Set-TextValue "D32" '8.01'
$ws.Range("E32").Value = '  -2.16%  '
Set-TextValue "D33" '2.12'
$ws.Range("E33").Value = '  -2.16%  '
$ws.Range("E34").Value = '  -0.06%  '
Set-TextValue "D35" '23.21'
$ws.Range("E35").Value = '  -0.68%  '
$ws.Range("E36").Value = '  -4.60%  '
Set-TextValue "D37" '1.53'
$ws.Range("E37").Value = '  -1.57%  '
Set-TextValue "D38" '6.76'
$ws.Range("E38").Value = '  -1.41%  '
Set-TextValue "D39" '164.70'
$ws.Range("E39").Value = '  -0.09%  '
Set-TextValue "D40" '0.0755'
$ws.Range("E40").Value = '  -3.80%  '
$ws.Range("E41").Value = '  +0.08%  '
$ws.Range("E42").Value = '  -1.09%  '
Set-TextValue "D43" '24.60'
$ws.Range("E43").Value = '  -1.05%  '
Set-TextValue "D44" '1.67'
$ws.Range("E44").Value = '  -3.05%  '
Set-TextValue "D45" '4.30'
$ws.Range("E45").Value = '  -2.45%  '
Set-TextValue "D46" '1.17'
$ws.Range("E46").Value = '  -5.28%  '
$ws.Range("D47").Value = '2.539.17'
$ws.Range("E47").Value = '  +8.70%  '
$ws.Range("B48").Value = 'InjectiveProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue "D48" '22.91'
$ws.Range("E48").Value = '  +0.66%  '
$ws.Range("B49").Value = 'Cosmos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue "D49" '6.75'
$ws.Range("E49").Value = '  -1.64%  '
Set-TextValue "D50" '2.41'
$ws.Range("E50").Value = '  +4.12%  '
$ws.Range("E51").Value = '  -1.53%  '
